$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 641.63635
$ws.Range("I2").Value = 605.9
$ws.Range("J2").Value = 999
$ws.Range("K2").Value = 605.9
$ws.Range("L2").Value = 999
$ws.Range("M2").Value = -492.9
$ws.Range("N2").Value = -1225

$ws.Range("H9").Value = 667001.3
$ws.Range("I9").Value = 420.27274
$ws.Range("K9").Value = 420.27274
$ws.Range("M9").Value = -251.27274

$ws.Range("H19").Value = 1225.5454
$ws.Range("J19").Value = 1275.7778
$ws.Range("L19").Value = 1275.7778
$ws.Range("N19").Value = -1625.7778

$ws.Range("H28").Value = 2692.8
$ws.Range("I28").Value = 2487.1667
$ws.Range("J28").Value = 3515.3333
$ws.Range("K28").Value = 2487.1667
$ws.Range("L28").Value = 3515.3333
$ws.Range("M28").Value = -2002.1667
$ws.Range("N28").Value = -4485.3333

$ws.Range("H41").Value = 727.8333
$ws.Range("I41").Value = 741.75
$ws.Range("K41").Value = 741.75
$ws.Range("M41").Value = -301.75

$ws.Range("H43").Value = 1450.8334
$ws.Range("I43").Value = 1461
$ws.Range("J43").Value = 1400
$ws.Range("K43").Value = 1461
$ws.Range("L43").Value = 1400
$ws.Range("M43").Value = -1392
$ws.Range("N43").Value = -1538

$ws.Range("H70").Value = 2576.5715
$ws.Range("I70").Value = 2423.25
$ws.Range("J70").Value = 2781
$ws.Range("K70").Value = 7269.75
$ws.Range("L70").Value = 8343
$ws.Range("M70").Value = -6999.75
$ws.Range("N70").Value = -8883

$ws.Range("H73").Value = 2576.5715
$ws.Range("I73").Value = 2423.25
$ws.Range("J73").Value = 2781
$ws.Range("K73").Value = 7269.75
$ws.Range("L73").Value = 8343
$ws.Range("M73").Value = -6333.75
$ws.Range("N73").Value = -10215

$ws.Range("H92").Value = 1294.6
$ws.Range("J92").Value = 5455
$ws.Range("L92").Value = 5455
$ws.Range("N92").Value = -7951

$ws.Range("H94").Value = 7130.727
$ws.Range("I94").Value = 6715.4443
$ws.Range("J94").Value = 8999.5
$ws.Range("K94").Value = 6715.4443
$ws.Range("L94").Value = 8999.5
$ws.Range("M94").Value = -6264.4443
$ws.Range("N94").Value = -9901.5

$ws.Range("H98").Value = 30667
$ws.Range("I98").Value = 34199.2
$ws.Range("J98").Value = 13006
$ws.Range("K98").Value = 34199.2
$ws.Range("L98").Value = 13006
$ws.Range("M98").Value = -32701.2
$ws.Range("N98").Value = -16002

$ws.Range("H112").Value = 14288026
$ws.Range("J112").Value = 28573492
$ws.Range("L112").Value = 85720476
$ws.Range("N112").Value = -85722692

$ws.Range("H122").Value = 30667
$ws.Range("I122").Value = 34199.2
$ws.Range("J122").Value = 13006
$ws.Range("K122").Value = 102597.6
$ws.Range("L122").Value = 39018
$ws.Range("M122").Value = -100147.6
$ws.Range("N122").Value = -43918

$ws.Range("H132").Value = 2449.3845
$ws.Range("I132").Value = 2327.9048
$ws.Range("K132").Value = 6983.714399999999
$ws.Range("M132").Value = -4453.714399999999

$ws.Range("H135").Value = 1877.1428
$ws.Range("I135").Value = 1856.8334
$ws.Range("K135").Value = 16711.5006
$ws.Range("M135").Value = -14176.5006

$ws.Range("H137").Value = 1161.8334
$ws.Range("I137").Value = 1171.375
$ws.Range("J137").Value = 1142.75
$ws.Range("K137").Value = 3514.125
$ws.Range("L137").Value = 3428.25
$ws.Range("M137").Value = -964.125
$ws.Range("N137").Value = -8528.25

$ws.Range("H138").Value = 687636.2
$ws.Range("J138").Value = 762467.4399999999
$ws.Range("L138").Value = 2287402.32
$ws.Range("N138").Value = -2297682.32

$ws.Range("H141").Value = 4091.2307
$ws.Range("I141").Value = 2289.6365
$ws.Range("J141").Value = 14000
$ws.Range("K141").Value = 6868.9095
$ws.Range("L141").Value = 42000
$ws.Range("M141").Value = -1688.9095
$ws.Range("N141").Value = -52360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1830.3158
$ws.Range("I2").Value = 1872.9375
$ws.Range("K2").Value = 1872.9375
$ws.Range("M2").Value = -1759.9375

$ws.Range("H3").Value = 5072.5
$ws.Range("I3").Value = 6716.6665
$ws.Range("J3").Value = 140
$ws.Range("K3").Value = 6716.6665
$ws.Range("L3").Value = 140
$ws.Range("M3").Value = -6601.6665
$ws.Range("N3").Value = -370

$ws.Range("H5").Value = 1293.8
$ws.Range("I5").Value = 1117.25
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 1117.25
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = -1005.25
$ws.Range("N5").Value = -2224

$ws.Range("H8").Value = 1402.8
$ws.Range("J8").Value = 199.5
$ws.Range("L8").Value = 199.5
$ws.Range("N8").Value = -487.5

$ws.Range("H11").Value = 9994
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 9994
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = None
$ws.Range("M11").Value = 9994
$ws.Range("N11").Value = -10282

$ws.Range("H12").Value = 248.4
$ws.Range("J12").Value = 497.5
$ws.Range("L12").Value = 497.5
$ws.Range("N12").Value = -843.5

$ws.Range("H32").Value = 2453.7805
$ws.Range("I32").Value = 1683.2703
$ws.Range("J32").Value = 9581
$ws.Range("K32").Value = 1683.2703
$ws.Range("L32").Value = 9581
$ws.Range("M32").Value = -1396.2703
$ws.Range("N32").Value = -10155

$ws.Range("H45").Value = 933.2857
$ws.Range("I45").Value = 855.9787
$ws.Range("J45").Value = 2750
$ws.Range("K45").Value = 855.9787
$ws.Range("L45").Value = 2750
$ws.Range("M45").Value = -478.9787
$ws.Range("N45").Value = -3504

$ws.Range("H61").Value = 2471.2727
$ws.Range("I61").Value = 2255.6191
$ws.Range("J61").Value = 7000
$ws.Range("K61").Value = 2255.6191
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = -2043.6191
$ws.Range("N61").Value = -7424

$ws.Range("H97").Value = 622
$ws.Range("I97").Value = 573.8333
$ws.Range("K97").Value = 573.8333
$ws.Range("M97").Value = -77.83330000000001

$ws.Range("H110").Value = 763.5833
$ws.Range("I110").Value = 666.5
$ws.Range("J110").Value = 1249
$ws.Range("K110").Value = 666.5
$ws.Range("L110").Value = 1249
$ws.Range("M110").Value = 1378.5
$ws.Range("N110").Value = -5339

$ws.Range("H116").Value = 1830.3158
$ws.Range("I116").Value = 1872.9375
$ws.Range("K116").Value = 1872.9375
$ws.Range("M116").Value = 421.0625

$ws.Range("H121").Value = 110000
$ws.Range("J121").Value = 110000
$ws.Range("L121").Value = 110000
$ws.Range("N121").Value = -113494

$ws.Range("H122").Value = 978.2
$ws.Range("I122").Value = 1130.3334
$ws.Range("J122").Value = 750
$ws.Range("K122").Value = 3391.0002
$ws.Range("L122").Value = 2250
$ws.Range("M122").Value = -941.0001999999999
$ws.Range("N122").Value = -7150

$ws.Range("H128").Value = 68749.5
$ws.Range("J128").Value = 68749.5
$ws.Range("L128").Value = 68749.5
$ws.Range("N128").Value = -78709.5

$ws.Range("H132").Value = 2746.0557
$ws.Range("I132").Value = 2251.3333
$ws.Range("K132").Value = 6753.999899999999
$ws.Range("M132").Value = -4223.999899999999

$ws.Range("H136").Value = 2471.2727
$ws.Range("I136").Value = 2255.6191
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 6766.8573
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -4216.8573
$ws.Range("N136").Value = -26100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1830.3158
$ws.Range("I3").Value = 1872.9375
$ws.Range("K3").Value = 1872.9375
$ws.Range("M3").Value = -1758.9375

$ws.Range("H4").Value = 1293.8
$ws.Range("I4").Value = 1117.25
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 1117.25
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = -1002.25
$ws.Range("N4").Value = -2230

$ws.Range("H5").Value = 2907.8572
$ws.Range("I5").Value = 2371
$ws.Range("J5").Value = 4250
$ws.Range("K5").Value = 2371
$ws.Range("L5").Value = 4250
$ws.Range("M5").Value = -2258
$ws.Range("N5").Value = -4476

$ws.Range("H20").Value = 1929.8096
$ws.Range("I20").Value = 1283.5
$ws.Range("K20").Value = 1283.5
$ws.Range("M20").Value = -1036.5

$ws.Range("H22").Value = 399.5
$ws.Range("I22").Value = 399.5
$ws.Range("K22").Value = 399.5
$ws.Range("M22").Value = -226.5

$ws.Range("H64").Value = 938.9
$ws.Range("I64").Value = 956.4286
$ws.Range("K64").Value = 956.4286
$ws.Range("M64").Value = -731.4286

$ws.Range("H67").Value = 938.9
$ws.Range("I67").Value = 956.4286
$ws.Range("K67").Value = 956.4286
$ws.Range("M67").Value = -176.4286

$ws.Range("H86").Value = 7284.5293
$ws.Range("I86").Value = 6742.533
$ws.Range("K86").Value = 6742.533
$ws.Range("M86").Value = -5619.533

$ws.Range("H89").Value = 7284.5293
$ws.Range("I89").Value = 6742.533
$ws.Range("K89").Value = 33712.665
$ws.Range("M89").Value = -28096.665

$ws.Range("H105").Value = 2391.1177
$ws.Range("I105").Value = 1874.6522
$ws.Range("J105").Value = 3471
$ws.Range("K105").Value = 1874.6522
$ws.Range("L105").Value = 3471
$ws.Range("M105").Value = -127.6522
$ws.Range("N105").Value = -6965

$ws.Range("H107").Value = 2714.2083
$ws.Range("I107").Value = 2708.5789
$ws.Range("J107").Value = 2735.6
$ws.Range("K107").Value = 2708.5789
$ws.Range("L107").Value = 2735.6
$ws.Range("M107").Value = -788.5789
$ws.Range("N107").Value = -6575.6

$ws.Range("H134").Value = 3753.6924
$ws.Range("I134").Value = 3687.6
$ws.Range("K134").Value = 11062.8
$ws.Range("M134").Value = -8527.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = None
$ws.Range("M2").Value = 1000
$ws.Range("N2").Value = -1226

$ws.Range("H6").Value = 7000799.5
$ws.Range("I6").Value = 17500526
$ws.Range("K6").Value = 17500526
$ws.Range("M6").Value = -17500413

$ws.Range("H7").Value = 681.1667
$ws.Range("I7").Value = 260.5
$ws.Range("K7").Value = 260.5
$ws.Range("M7").Value = -147.5

$ws.Range("H22").Value = 1697.6666
$ws.Range("I22").Value = 1547
$ws.Range("J22").Value = 1999
$ws.Range("K22").Value = 1547
$ws.Range("L22").Value = 1999
$ws.Range("M22").Value = -1197
$ws.Range("N22").Value = -2699

$ws.Range("H31").Value = 2147.2104
$ws.Range("I31").Value = 1890.4286
$ws.Range("J31").Value = 2866.2
$ws.Range("K31").Value = 1890.4286
$ws.Range("L31").Value = 2866.2
$ws.Range("M31").Value = -1595.4286
$ws.Range("N31").Value = -3456.2

$ws.Range("H34").Value = 2147.2104
$ws.Range("I34").Value = 1890.4286
$ws.Range("J34").Value = 2866.2
$ws.Range("K34").Value = 1890.4286
$ws.Range("L34").Value = 2866.2
$ws.Range("M34").Value = -1688.4286
$ws.Range("N34").Value = -3270.2

$ws.Range("H41").Value = 19665.666
$ws.Range("J41").Value = 29999
$ws.Range("L41").Value = 29999
$ws.Range("N41").Value = -30855

$ws.Range("H80").Value = 32998.332
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 32998.332
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = None
$ws.Range("M80").Value = 32998.332
$ws.Range("N80").Value = -35244.332

$ws.Range("H83").Value = 32998.332
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 32998.332
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = None
$ws.Range("M83").Value = 98994.99600000001
$ws.Range("N83").Value = -110226.996

$ws.Range("H86").Value = 6252.75
$ws.Range("I86").Value = 6337.3335
$ws.Range("J86").Value = 5999
$ws.Range("K86").Value = 6337.3335
$ws.Range("L86").Value = 5999
$ws.Range("M86").Value = -5214.3335
$ws.Range("N86").Value = -8245

$ws.Range("H89").Value = 6252.75
$ws.Range("I89").Value = 6337.3335
$ws.Range("J89").Value = 5999
$ws.Range("K89").Value = 31686.6675
$ws.Range("L89").Value = 29995
$ws.Range("M89").Value = -26070.6675
$ws.Range("N89").Value = -41227

$ws.Range("H99").Value = 24248.875
$ws.Range("I99").Value = 35939.8
$ws.Range("J99").Value = 18934.818
$ws.Range("K99").Value = 35939.8
$ws.Range("L99").Value = 18934.818
$ws.Range("M99").Value = -34441.8
$ws.Range("N99").Value = -21930.818

$ws.Range("H126").Value = 24248.875
$ws.Range("I126").Value = 35939.8
$ws.Range("J126").Value = 18934.818
$ws.Range("K126").Value = 107819.4
$ws.Range("L126").Value = 56804.454
$ws.Range("M126").Value = -105349.4
$ws.Range("N126").Value = -61744.454

$ws.Range("H132").Value = 3206.4119
$ws.Range("I132").Value = 3191.2144
$ws.Range("K132").Value = 9573.643199999999
$ws.Range("M132").Value = -7043.643199999999

$ws.Range("H134").Value = 4373.467
$ws.Range("I134").Value = 4638.6924
$ws.Range("J134").Value = 2649.5
$ws.Range("K134").Value = 13916.0772
$ws.Range("L134").Value = 7948.5
$ws.Range("M134").Value = -11381.0772
$ws.Range("N134").Value = -13018.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1648
$ws.Range("J2").Value = 20
$ws.Range("L2").Value = 120
$ws.Range("N2").Value = -346

$ws.Range("H17").Value = 616.6667
$ws.Range("J17").Value = 800
$ws.Range("L17").Value = 2400
$ws.Range("N17").Value = -2738

$ws.Range("H26").Value = 362.3125
$ws.Range("I26").Value = 34.4
$ws.Range("K26").Value = 103.2
$ws.Range("M26").Value = 184.8

$ws.Range("H33").Value = 604.6667
$ws.Range("I33").Value = 323.66666
$ws.Range("J33").Value = 1166.6666
$ws.Range("K33").Value = 1941.99996
$ws.Range("L33").Value = 6999.9996
$ws.Range("M33").Value = -1658.99996
$ws.Range("N33").Value = -7565.9996

$ws.Range("H38").Value = 378.45456
$ws.Range("I38").Value = 23.75
$ws.Range("J38").Value = 581.1429000000001
$ws.Range("K38").Value = 71.25
$ws.Range("L38").Value = 1743.4287
$ws.Range("M38").Value = 275.75
$ws.Range("N38").Value = -2437.4287

$ws.Range("H54").Value = 2200
$ws.Range("J54").Value = 3000
$ws.Range("L54").Value = 9000
$ws.Range("N54").Value = -10118

$ws.Range("H68").Value = 2022.127
$ws.Range("J68").Value = 2270.8823
$ws.Range("L68").Value = 6812.646900000001
$ws.Range("N68").Value = -8434.6469

$ws.Range("H71").Value = 2022.127
$ws.Range("J71").Value = 2270.8823
$ws.Range("L71").Value = 20437.9407
$ws.Range("N71").Value = -28549.9407

$ws.Range("H86").Value = 639.7273
$ws.Range("J86").Value = 802.6
$ws.Range("L86").Value = 2407.8
$ws.Range("N86").Value = -4779.8

$ws.Range("H89").Value = 639.7273
$ws.Range("J89").Value = 802.6
$ws.Range("L89").Value = 7223.400000000001
$ws.Range("N89").Value = -19079.4

$ws.Range("H110").Value = 19981
$ws.Range("I110").Value = 9962
$ws.Range("K110").Value = 29886
$ws.Range("M110").Value = -25796

$ws.Range("H117").Value = 3487.0908
$ws.Range("J117").Value = 4084.4443
$ws.Range("L117").Value = 12253.3329
$ws.Range("N117").Value = -19137.3329

$ws.Range("H131").Value = 17404.215
$ws.Range("I131").Value = 112230.8
$ws.Range("J131").Value = 1599.7833
$ws.Range("K131").Value = 336692.4
$ws.Range("L131").Value = 4799.3499
$ws.Range("M131").Value = -331652.4
$ws.Range("N131").Value = -14879.3499

$ws.Range("H132").Value = 113.75
$ws.Range("J132").Value = 155
$ws.Range("L132").Value = 1395
$ws.Range("N132").Value = -6455

$ws.Range("H134").Value = 1629.4445
$ws.Range("I134").Value = 1629.4445
$ws.Range("K134").Value = 4888.333500000001
$ws.Range("M134").Value = 181.6664999999994

$ws.Range("H137").Value = 17679.125
$ws.Range("J137").Value = 17777.572
$ws.Range("L137").Value = 53332.716
$ws.Range("N137").Value = -63532.716

$ws.Range("H138").Value = 1488.1666
$ws.Range("I138").Value = 1488.1666
$ws.Range("K138").Value = 4464.4998
$ws.Range("M138").Value = 675.5002000000004

$ws.Range("H139").Value = 8574.223
$ws.Range("I139").Value = 6879.8335
$ws.Range("J139").Value = 11963
$ws.Range("K139").Value = 20639.5005
$ws.Range("L139").Value = 35889
$ws.Range("M139").Value = -15499.5005
$ws.Range("N139").Value = -46169

$ws.Range("H140").Value = 1625.2632
$ws.Range("I140").Value = 1625.2632
$ws.Range("K140").Value = 4875.7896
$ws.Range("M140").Value = 304.2103999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 328.64285
$ws.Range("I2").Value = 328.64285
$ws.Range("K2").Value = 328.64285
$ws.Range("M2").Value = -215.64285

$ws.Range("H74").Value = 90000
$ws.Range("J74").Value = 90000
$ws.Range("L74").Value = 90000
$ws.Range("N74").Value = -91872

$ws.Range("H77").Value = 90000
$ws.Range("J77").Value = 90000
$ws.Range("L77").Value = 270000
$ws.Range("N77").Value = -279360

$ws.Range("H80").Value = 4642.875
$ws.Range("I80").Value = 3027.7
$ws.Range("K80").Value = 3027.7
$ws.Range("M80").Value = -2029.7

$ws.Range("H83").Value = 4642.875
$ws.Range("I83").Value = 3027.7
$ws.Range("K83").Value = 15138.5
$ws.Range("M83").Value = -10146.5

$ws.Range("H102").Value = 4399.5454
$ws.Range("I102").Value = 4395.1
$ws.Range("J102").Value = 4444
$ws.Range("K102").Value = 4395.1
$ws.Range("L102").Value = 4444
$ws.Range("M102").Value = -2773.1
$ws.Range("N102").Value = -7688

$ws.Range("H107").Value = 447.44446
$ws.Range("I107").Value = 259.15384
$ws.Range("K107").Value = 259.15384
$ws.Range("M107").Value = 1660.84616

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = None
$ws.Range("N117").Value = 0

$ws.Range("H122").Value = 1349.3334
$ws.Range("I122").Value = 1349.3334
$ws.Range("K122").Value = 4048.0002
$ws.Range("M122").Value = -1598.0002

$ws.Range("H132").Value = 3951.1538
$ws.Range("I132").Value = 3255.875
$ws.Range("K132").Value = 9767.625
$ws.Range("M132").Value = -7237.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3115.111
$ws.Range("I7").Value = 2576.7144
$ws.Range("K7").Value = 2576.7144
$ws.Range("M7").Value = -2464.7144

$ws.Range("H22").Value = 1674.75
$ws.Range("I22").Value = 1533.1666
$ws.Range("J22").Value = 2099.5
$ws.Range("K22").Value = 1533.1666
$ws.Range("L22").Value = 2099.5
$ws.Range("M22").Value = -1238.1666
$ws.Range("N22").Value = -2689.5

$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = None
$ws.Range("M26").Value = None
$ws.Range("N26").Value = 0

$ws.Range("H27").Value = 1674.75
$ws.Range("I27").Value = 1533.1666
$ws.Range("J27").Value = 2099.5
$ws.Range("K27").Value = 1533.1666
$ws.Range("L27").Value = 2099.5
$ws.Range("M27").Value = -1426.1666
$ws.Range("N27").Value = -2313.5

$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = None
$ws.Range("N29").Value = 0

$ws.Range("H40").Value = 2085.125
$ws.Range("I40").Value = 1536.8
$ws.Range("J40").Value = 2999
$ws.Range("K40").Value = 1536.8
$ws.Range("L40").Value = 2999
$ws.Range("M40").Value = -1400.8
$ws.Range("N40").Value = -3271

$ws.Range("H55").Value = 198.90909
$ws.Range("I55").Value = 187.33333
$ws.Range("K55").Value = 187.33333
$ws.Range("M55").Value = -14.33332999999999

$ws.Range("H61").Value = 7433.1113
$ws.Range("I61").Value = 7160.2144
$ws.Range("K61").Value = 7160.2144
$ws.Range("M61").Value = -6958.2144

$ws.Range("H113").Value = 7433.1113
$ws.Range("I113").Value = 7160.2144
$ws.Range("K113").Value = 7160.2144
$ws.Range("M113").Value = -4990.2144

$ws.Range("H122").Value = 2778
$ws.Range("I122").Value = 2778
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8334
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = None
$ws.Range("N122").Value = -5884

$ws.Range("H126").Value = 3115.111
$ws.Range("I126").Value = 2576.7144
$ws.Range("K126").Value = 7730.1432
$ws.Range("M126").Value = -5260.1432

$ws.Range("H132").Value = 5174.5557
$ws.Range("I132").Value = 4403.2144
$ws.Range("K132").Value = 13209.6432
$ws.Range("M132").Value = -10679.6432

$ws.Range("H136").Value = 3494.476
$ws.Range("I136").Value = 3346.3572
$ws.Range("J136").Value = 3790.7144
$ws.Range("K136").Value = 10039.0716
$ws.Range("L136").Value = 11372.1432
$ws.Range("M136").Value = -7489.071599999999
$ws.Range("N136").Value = -16472.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = None
$ws.Range("N9").Value = 0

$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = None
$ws.Range("M32").Value = None
$ws.Range("N32").Value = 0

$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").Value = None

$ws.Range("H62").Value = 9975.4
$ws.Range("J62").Value = 9975.4
$ws.Range("L62").Value = 9975.4
$ws.Range("N62").Value = -11223.4

$ws.Range("H65").Value = 9975.4
$ws.Range("J65").Value = 9975.4
$ws.Range("L65").Value = 49877
$ws.Range("N65").Value = -56117

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = None
$ws.Range("N87").Value = 0

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = None
$ws.Range("N90").Value = 0

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = None
$ws.Range("N92").Value = 0

$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = None
$ws.Range("N94").Value = 0

$ws.Range("H107").Value = 1391.5
$ws.Range("J107").Value = 1499.6666
$ws.Range("L107").Value = 4498.9998
$ws.Range("N107").Value = -8338.9998

$ws.Range("H110").Value = 141966.33
$ws.Range("J110").Value = 141966.33
$ws.Range("L110").Value = 141966.33
$ws.Range("N110").Value = -150146.33

$ws.Range("H122").Value = 3683.4
$ws.Range("I122").Value = 4031.7368
$ws.Range("K122").Value = 12095.2104
$ws.Range("M122").Value = -9645.2104

$ws.Range("H126").Value = 2424.5
$ws.Range("I126").Value = 2424.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7273.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = None
$ws.Range("N126").Value = -4803.5

$ws.Range("H132").Value = 4884.452
$ws.Range("I132").Value = 3436.5151
$ws.Range("J132").Value = 10193.556
$ws.Range("K132").Value = 10309.5453
$ws.Range("L132").Value = 30580.668
$ws.Range("M132").Value = -7779.5453
$ws.Range("N132").Value = -35640.66800000001

$ws.Range("H136").Value = 3240.484
$ws.Range("I136").Value = 3824.182
$ws.Range("K136").Value = 11472.546
$ws.Range("M136").Value = -8922.545999999998
